$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from A67 down to A68:A71 so new date cells inherit the date number format
$ws.Range("A67").Copy()
$ws.Range("A68:A71").PasteSpecial(-4122)  # xlPasteFormats

# Row 68: 2018-12-07, Daniel, styling changes, 4
$ws.Range("A68").Value = 43441
$ws.Range("B68").Value = "Daniel"
$ws.Range("C68").Value = "styling changes"
$ws.Range("D68").Value = 4

# Row 69: 2018-12-10, Daniel, styling changes, 3
$ws.Range("A69").Value = 43444
$ws.Range("B69").Value = "Daniel"
$ws.Range("C69").Value = "styling changes"
$ws.Range("D69").Value = 3

# Row 70: 2018-12-11, Daniel, styling changes, 3
$ws.Range("A70").Value = 43445
$ws.Range("B70").Value = "Daniel"
$ws.Range("C70").Value = "styling changes"
$ws.Range("D70").Value = 3

# Row 71: empty date cell with style only
$ws.Range("F66").Select()
